$d = $word.ActiveDocument

# Paragraph 1 (originally the sole empty paragraph) becomes the title.
$p1 = $d.Paragraphs(1)
$p1.Range.Text = "Introduction to Version Control in GitHub"
$p1.Range.InsertParagraphAfter()

# Paragraph 2 is the original empty paragraph, now followed by a new
# paragraph of body text.
$p2 = $d.Paragraphs(2)
$p2.Range.InsertParagraphAfter()

$p3 = $d.Paragraphs(3)
$p3.Range.Text = "Version control is a system that helps developers manage and track changes to their code over time. It allows multiple contributors to work on the same project without conflicts, providing a history of all modifications and making collaboration efficient and reliable."
$p3.Range.InsertParagraphAfter()

$p4 = $d.Paragraphs(4)
$p4.Range.InsertParagraphAfter()

$p5 = $d.Paragraphs(5)
$p5.Range.Text = "GitHub is a web-based platform built around Git, the world’s most widely used distributed version control system. It allows individuals and teams to host, share, and collaborate on code repositories. With GitHub, developers can manage versions of their code, collaborate through pull requests, review changes, and integrate automated workflows to build and deploy software."
$p5.Range.InsertParagraphAfter()

$p6 = $d.Paragraphs(6)
$p6.Range.InsertParagraphAfter()

$p7 = $d.Paragraphs(7)
$p7.Range.Text = "GitHub is widely used by both open-source communities and enterprises due to its simplicity, flexibility, and integration with various development tools. It provides not only version control but also a complete ecosystem for project management, automation, and collaboration."
